# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Re-orders / updates the three "employee" rows (17-19) on Hoja1 so the
# account-statement data reflects the refreshed database export: the
# "OBELIS EDMUNDO PEÑA SUAREZ" record now appears first (row 17), the
# "LUIS DAVID CEBALLOS PEREZ" record appears twice (rows 18-19, for periods
# 1711 and 1710 respectively), and the "Valor Mora" / "Salario Basico"
# amounts are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 17: now OBELIS EDMUNDO PEÑA SUAREZ, periodo 1711
$ws.Range("C17").Value = "9296620"
$ws.Range("D17").Value = "OBELIS EDMUNDO PEÑA SUAREZ"
$ws.Range("E17").Value = "1711"
$ws.Range("F17").Value = 24369
$ws.Range("G17").Value = 1218448

# Row 18: LUIS DAVID CEBALLOS PEREZ, periodo 1711
$ws.Range("C18").Value = "1063170757"
$ws.Range("D18").Value = "LUIS DAVID CEBALLOS PEREZ"
$ws.Range("E18").Value = "1711"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 1206660

# Row 19: LUIS DAVID CEBALLOS PEREZ, periodo 1710
$ws.Range("C19").Value = "1063170757"
$ws.Range("D19").Value = "LUIS DAVID CEBALLOS PEREZ"
$ws.Range("E19").Value = "1710"
$ws.Range("F19").Value = 14755
$ws.Range("G19").Value = 1206660
